# Update with Correct Forecast output
# - Renames Sheet1 to "Sales vs PO"
# - Adds three new sheets: "Weekly Growth", "Volume Insights", "Prediction Info"
# - Inserts a new "Order Week" column (old ds values) into the first sheet and
#   shifts the old PO_Requested_Qty values out (now zeroed there)
# - Populates the new sheets with the weekly growth, summary stats and forecast

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- create the remaining sheets, in order, right after the first one -------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

$ws1.Name = "Sales vs PO"

# Grab the bold/centered header style that's already used on row 1 of the
# first sheet so the new sheets' headers match without minting new styles.
$headerStyleSource = $ws1.Range("A1")

# =============================================================================
# Sheet 1 : "Sales vs PO" -- insert an "Order Week" column before the old
# PO_Requested_Qty column (which shifts right, and is zeroed out since the
# real PO numbers now live on the "Weekly Growth" sheet).
# =============================================================================
$ws1.Columns.Item(3).Insert()

$ws1.Range("C1").Value = "Order Week"
$ws1.Range("C2:C19").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$newDs = @(45536,45543,45550,45557,45564,45571,45578,45585,45592,45599,45606,45613,45620,45627,45634,45641,45648,45655)
$oldDs = @(45530,45537,45544,45551,45558,45565,45572,45579,45586,45593,45600,45607,45614,45621,45628,45635,45642,45649)

for ($i = 0; $i -lt $newDs.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 1).Value = $newDs[$i]
    $ws1.Cells.Item($row, 3).Value = $oldDs[$i]
    $ws1.Cells.Item($row, 4).Value = 0
}

# =============================================================================
# Sheet 2 : "Weekly Growth" -- the non-zero weekly PO quantities plus their
# week-over-week growth percentage.
# =============================================================================
$ws2.Range("A1").Value = "ds"
$ws2.Range("B1").Value = "PO_Requested_Qty"
$ws2.Range("C1").Value = "Growth%"

$gDs     = @(45537, 45558, 45586)
$gQty    = @(32, 16, 16)
$gGrowth = @(0, -50, 0)

for ($i = 0; $i -lt $gDs.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $gDs[$i]
    $ws2.Cells.Item($row, 2).Value = $gQty[$i]
    $ws2.Cells.Item($row, 3).Value = $gGrowth[$i]
}
$ws2.Range("A2:A4").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$headerStyleSource.Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)

# =============================================================================
# Sheet 3 : "Volume Insights" -- summary statistics over the PO quantities.
# =============================================================================
$ws3.Range("A1").Value = "Total_PO_Quantity"
$ws3.Range("B1").Value = "Average_PO_Quantity"
$ws3.Range("C1").Value = "Max_PO_Quantity"
$ws3.Range("D1").Value = "Min_PO_Quantity"

$ws3.Range("A2").Value = 64
$ws3.Range("B2").Value = 21.33333333333333
$ws3.Range("C2").Value = 32
$ws3.Range("D2").Value = 16

$headerStyleSource.Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# =============================================================================
# Sheet 4 : "Prediction Info" -- the forecast for next week's PO quantity.
# =============================================================================
$ws4.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"
$ws4.Range("A2").Value = 5.333333333333343

$headerStyleSource.Copy()
$ws4.Range("A1").PasteSpecial(-4122)

$ws1.Select()
